# Applies the textual edits from the commit "minor edits, close #30"
# to the Response to Reviewers document. Each block performs a single
# Find & Replace (wildcards off, plain text) over the whole document body.
$d = $word.ActiveDocument

$failures = 0

# 1. remove 'additional'
$rng = $d.Content
$ok = $rng.Find.Execute("conducted additional sensitivity analyses to show the robustness of our results", $true, $false, $false, $false, $false, $true, 1, $false, "conducted sensitivity analyses to show the robustness of our results", 2)
if (-not $ok) {
    Write-Host "FAILED: remove 'additional'"
    $failures = $failures + 1
}

# 2. rewrite conceptual clarifications sentence
$rng = $d.Content
$ok = $rng.Find.Execute(" made some conceptual clarifications, including our contribution to the literature as a theoretical work and a clarification that we took sex a main explanatory variable for mating preferences and investigated its different impact on behaviors among heterosexual individuals and gay males and lesbian females.", $true, $false, $false, $false, $false, $true, 1, $false, " made some conceptual clarifications, including our contribution to the literature as a piece of theoretical work and our focus on sex as a main explanatory variable for mating preferences and its impact on behaviors among different populations.", 2)
if (-not $ok) {
    Write-Host "FAILED: rewrite conceptual clarifications sentence"
    $failures = $failures + 1
}

# 3. orientations -> orientation #1
$rng = $d.Content
$ok = $rng.Find.Execute("regardless of sexual orientations (we already fully explicated the reasoning and evidence for this claim earlier in the Introduction).", $true, $false, $false, $false, $false, $true, 1, $false, "regardless of sexual orientation (we already fully explicated the reasoning and evidence for this claim earlier in the Introduction).", 2)
if (-not $ok) {
    Write-Host "FAILED: orientations -> orientation #1"
    $failures = $failures + 1
}

# 4. orientations -> orientation #2 + add among males
$rng = $d.Content
$ok = $rng.Find.Execute("regardless of sexual orientations, this interest can lead to different behavioral consequences depending on whether they are heterosexual (partners being", $true, $false, $false, $false, $false, $true, 1, $false, "regardless of sexual orientation, this interest can lead to different behavioral consequences among males depending on whether they are heterosexual (partners being", 2)
if (-not $ok) {
    Write-Host "FAILED: orientations -> orientation #2 + add among males"
    $failures = $failures + 1
}

# 5. github repo: when -> after the
$rng = $d.Content
$ok = $rng.Find.Execute(". The Github repository has been temporarily unlinked but will be linked back when this paper gets accepted.", $true, $false, $false, $false, $false, $true, 1, $false, ". The Github repository has been temporarily unlinked but will be linked back after the paper gets accepted.", 2)
if (-not $ok) {
    Write-Host "FAILED: github repo: when -> after the"
    $failures = $failures + 1
}

# 6. sensitivity analyses results reorder
$rng = $d.Content
$ok = $rng.Find.Execute("The results from the sensitivity analyses were similar to the main results presented in the manuscript and were reported in the Supplemental Materials. In the revised manuscript,", $true, $false, $false, $false, $false, $true, 1, $false, "The results from the sensitivity analyses, now reported in the Supplemental Materials, were similar to the main results presented in the manuscript. In the revised manuscript,", 2)
if (-not $ok) {
    Write-Host "FAILED: sensitivity analyses results reorder"
    $failures = $failures + 1
}

# 7. add 'of the Supplemental Materials'
$rng = $d.Content
$ok = $rng.Find.Execute("he likelihood of short-term mating for paired agents was reduced (See the “Process overview and scheduling” for a detailed description of the procedures).", $true, $false, $false, $false, $false, $true, 1, $false, "he likelihood of short-term mating for paired agents was reduced (See the “Process overview and scheduling” of the Supplemental Materials for a detailed description of the procedures).", 2)
if (-not $ok) {
    Write-Host "FAILED: add 'of the Supplemental Materials'"
    $failures = $failures + 1
}

# 8. rephrase engaging sentence
$rng = $d.Content
$ok = $rng.Find.Execute("engaging in short-term mating is possible while staying in a long-term relationship", $true, $false, $false, $false, $false, $true, 1, $false, "it is possible to engage in short-term mating while staying in a long-term relationship", 2)
if (-not $ok) {
    Write-Host "FAILED: rephrase engaging sentence"
    $failures = $failures + 1
}

# 9. rewrite agree with reviewer sentence
$rng = $d.Content
$ok = $rng.Find.Execute("Response: We agree with the reviewer on this point. The direct reason for the unrealistically large effect sizes is that there", $true, $false, $false, $false, $false, $true, 1, $false, "Response: We agree with the reviewer that the effect sizes are extremely large. We think the direct reason is that there", 2)
if (-not $ok) {
    Write-Host "FAILED: rewrite agree with reviewer sentence"
    $failures = $failures + 1
}

# 10. of -> between group means
$rng = $d.Content
$ok = $rng.Find.Execute("the differences of the group means. To avoid confusion, we reported", $true, $false, $false, $false, $false, $true, 1, $false, "the differences between the group means. To avoid confusion, we reported", 2)
if (-not $ok) {
    Write-Host "FAILED: of -> between group means"
    $failures = $failures + 1
}

# 11. remove 'further'
$rng = $d.Content
$ok = $rng.Find.Execute("in the revised manuscript to avoid further confusion.", $true, $false, $false, $false, $false, $true, 1, $false, "in the revised manuscript to avoid confusion.", 2)
if (-not $ok) {
    Write-Host "FAILED: remove 'further'"
    $failures = $failures + 1
}

# 12. i.e. -> that is
$rng = $d.Content
$ok = $rng.Find.Execute("i.e., in mating behaviors of short duration without commitment (Buss & Schmitt, 1993)", $true, $false, $false, $false, $false, $true, 1, $false, "that is, in mating behaviors of short duration without commitment (Buss & Schmitt, 1993)", 2)
if (-not $ok) {
    Write-Host "FAILED: i.e. -> that is"
    $failures = $failures + 1
}

# 13. agent-based modeling paragraph rewrite
$rng = $d.Content
$ok = $rng.Find.Execute("Response: Our manuscript reports a study based on agent-based modeling, which is a theoretical contribution rather than an empirical one. The study did not involve any participants.", $true, $false, $false, $false, $false, $true, 1, $false, "Response: Our manuscript reports a study based on agent-based modeling which did not involve any participants. Therefore, we are not able to provide demographic data about the participants or the method of selecting participants.", 2)
if (-not $ok) {
    Write-Host "FAILED: agent-based modeling paragraph rewrite"
    $failures = $failures + 1
}

# 14. worth investigating empirically -> worth empirical investigations
$rng = $d.Content
$ok = $rng.Find.Execute("issues that are worth investigating empirically. For example,", $true, $false, $false, $false, $false, $true, 1, $false, "issues that are worth empirical investigations. For example,", 2)
if (-not $ok) {
    Write-Host "FAILED: worth investigating empirically -> worth empirical investigations"
    $failures = $failures + 1
}

# 15. conceptual link -> causal relation
$rng = $d.Content
$ok = $rng.Find.Execute("a) the conceptual link between", $true, $false, $false, $false, $false, $true, 1, $false, "a) the causal relation between", 2)
if (-not $ok) {
    Write-Host "FAILED: conceptual link -> causal relation"
    $failures = $failures + 1
}

# 16. add 'mating' before behaviors
$rng = $d.Content
$ok = $rng.Find.Execute("differences in mating preferences and those in behaviors, and b) the existence of a shared explanation for", $true, $false, $false, $false, $false, $true, 1, $false, "differences in mating preferences and those in mating behaviors, and b) the existence of a shared explanation for", 2)
if (-not $ok) {
    Write-Host "FAILED: add 'mating' before behaviors"
    $failures = $failures + 1
}

# 17. gender/sex wording paragraph rewrite
$rng = $d.Content
$ok = $rng.Find.Execute("Response: In light of the broad range of evidence, we have changed the wording throughout the edited manuscript. Specifically, “gender difference” was replaced by “sex difference”, “men” was replaced by “males”, and “women” was replaced by “females” (except for discussions of existing studies as they mostly use the notion of “gender” rather than “sex).", $true, $false, $false, $false, $false, $true, 1, $false, "Response: In light of the broad range of evidence, we have changed the wording throughout the revised manuscript where it is appropriate. Specifically, “gender difference” was replaced by “sex difference”, “men” was replaced by “males”, and “women” was replaced by “females”.", 2)
if (-not $ok) {
    Write-Host "FAILED: gender/sex wording paragraph rewrite"
    $failures = $failures + 1
}

if ($failures -eq 0) {
    Write-Host "All replacements applied successfully."
} else {
    Write-Host "$failures replacement(s) failed."
}
